$d = $word.ActiveDocument

# --- Typo fixes -----------------------------------------------------------

# Fix "ride-lailing" -> "ride-hailing"
$d.Content.Find.Execute("ride-lailing", $false, $false, $false, $false, $false, $true, 1, $false, "ride-hailing", 2) | Out-Null

# Fix "panls" -> "plans"
$d.Content.Find.Execute("panls", $false, $false, $false, $false, $false, $true, 1, $false, "plans", 2) | Out-Null

# Fix "less like ly" -> "less likely"
$d.Content.Find.Execute("less like ly", $false, $false, $false, $false, $false, $true, 1, $false, "less likely", 2) | Out-Null

# --- Merge the split bold vocabulary words --------------------------------
# These words were previously split into two runs (e.g. "A" + "utomotive").
# Replacing the whole word re-merges it into a single run, but the engine
# also folds the immediately-following ": " run into it (same formatting).
# Re-toggling Bold on just the ": " restores it as its own run again.

$words = @("Automotive", "Unevenly", "Scratch", "Oversee", "Progressed")
foreach ($w in $words) {
    $d.Content.Find.Execute($w, $true, $false, $false, $false, $false, $true, 1, $false, $w, 2) | Out-Null

    $rng = $d.Content
    $rng.Find.Execute($w, $true) | Out-Null
    $colon = $d.Range($rng.End, $rng.End + 2)
    if ($colon.Text -eq ": ") {
        $colon.Font.Bold = $false
        $colon.Font.Bold = $true
    }
}
